$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "No"
$ws.Range("C9").Value = "Yes"
$ws.Range("C51").Value = "No"
$ws.Range("C53").Value = "Yes"
$ws.Range("C55").Value = "Yes"
$ws.Range("C56").Value = "No"
$ws.Range("C64").Value = "Yes"
$ws.Range("C66").Value = "Yes"
$ws.Range("C69").Value = "No"
$ws.Range("C71").Value = "No"
$ws.Range("C75").Value = "Yes"
$ws.Range("C78").Value = "No"
